$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$ws.Range("A2").Value = "0.27478192589996686 - 0.185445782038399y_1 + 0.38865877574366636y_2"
$ws.Range("B2").Value = "-0.27478192589996686"
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").Value = "0.03"
$ws.Range("E2").Value = "9.9"
$ws.Range("F2").Value = "0"

$ws.Range("A3").Value = "-10.578432249307559 + 1.3716822688057193y_1 + 0.17734556199482976y_2"
$ws.Range("B3").Value = "6.578432249307559"
$ws.Range("C3").Value = "J_0_L0_v"
$ws.Range("D3").Value = "0.85"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "3.5999999999999996"

$ws.Range("A4").Value = "52.61516996626992 - 2x - 5.489197429859824y_1 + 0.9037254082856405y_2"
$ws.Range("B4").Value = "-68.61516996626992"
$ws.Range("C4").Value = "J_0_LP_v"
$ws.Range("D4").Value = "0.73"
$ws.Range("E4").Value = "9.200000000000001"
$ws.Range("F4").Value = "9.399999999999999"

$ws.Range("A5").Value = "-87.49497027933667 + 8x + 3.315792066563272y_1 + 1.104963782177365y_2"
$ws.Range("B5").Value = "38.66497027933667"
$ws.Range("C5").Value = "J_Ne_L0_v"
$ws.Range("D5").Value = "0.02"
$ws.Range("E5").Value = "9.8"
$ws.Range("F5").Value = "9.5"

$ws.Range("A6").Value = "4.555913528392346 - 2x + 0.8153449226905218y_1 + 1.3433218892258285y_2"
$ws.Range("B6").Value = "6.954086471607654"
$ws.Range("C6").Value = "J_Ne_L0_v"
$ws.Range("D6").Value = "0.72"
$ws.Range("E6").Value = "0"
$ws.Range("F6").Value = "3.5999999999999996"

# --- Sheet: Punto_modificado ---
$ws2 = $wb.Worksheets.Item("Punto_modificado")
$ws2.Range("A2").Value = "7.4"
$ws2.Range("B2").Value = "7.35"
$ws2.Range("C2").Value = "2.8"

# --- Sheet: Vector_bf ---
$ws3 = $wb.Worksheets.Item("Vector_bf")
$ws3.Range("A2").Value = "3.1933833831055205"
$ws3.Range("A3").Value = "-1.8114140749025767"

# --- Sheet: Vector_BF ---
$ws4 = $wb.Worksheets.Item("Vector_BF")
$ws4.Range("A2").Value = "-59.0"
$ws4.Range("A3").Value = "22.84176734457047"
$ws4.Range("A4").Value = "-24.99064070142837"

# --- Sheet: Vector_Alpha ---
$ws5 = $wb.Worksheets.Item("Vector_Alpha")
$ws5.Range("A2").Value = 2.472015670882792
$ws5.Range("A3").Value = 1.1795047684717235
